$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: transaction / item row gets populated with data
$ws.Range("A4").Value = 1

# B4:G4 and N4 share one style class; switch it (and H4:K4's class) to
# Text format so the values below aren't re-interpreted as time/number.
$ws.Range("B4:G4").NumberFormat = "@"
$ws.Range("N4").NumberFormat = "@"
$ws.Range("H4:K4").NumberFormat = "@"

$ws.Range("B4").Value = "1 2 3 (ONE TWO THREE) 20 F.C.TABS."
$ws.Range("H4").Value = "1:0"
$ws.Range("L4").Value = 80
$ws.Range("N4").Value = "2:0"

# Row 5: totals row grows slightly taller and gets an amount
$ws.Rows.Item(5).RowHeight = 26.25
$ws.Range("K5").Value = 80
